$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.135036496350365
$ws.Range("C2").Value = 0.6642335766423357
$ws.Range("J2").Value = 0.0145985401459854
$ws.Range("P2").Value = 0.1021897810218978
$ws.Range("S2").Value = 0.08394160583941605
$ws.Range("B3").Value = 0.0053475935828877
$ws.Range("C3").Value = 0.0106951871657754
$ws.Range("J3").Value = 0.0267379679144385
$ws.Range("P3").Value = 0.7593582887700535
$ws.Range("S3").Value = 0.1978609625668449
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.04950495049504951
$ws.Range("D6").Value = 0.0198019801980198
$ws.Range("F6").Value = 0.05445544554455446
$ws.Range("J6").Value = 0.2722772277227723
$ws.Range("O6").Value = 0.004950495049504951
$ws.Range("Q6").Value = 0.2128712871287129
$ws.Range("R6").Value = 0.0594059405940594
$ws.Range("S6").Value = 0.3267326732673267
$ws.Range("B7").Value = 0.145945945945946
$ws.Range("D7").Value = 0.02702702702702703
$ws.Range("E7").Value = 0.01081081081081081
$ws.Range("F7").Value = 0.05945945945945946
$ws.Range("J7").Value = 0.06486486486486487
$ws.Range("O7").Value = 0.005405405405405406
$ws.Range("Q7").Value = 0.227027027027027
$ws.Range("R7").Value = 0.0918918918918919
$ws.Range("S7").Value = 0.3675675675675676
$ws.Range("B8").Value = 0.09033613445378151
$ws.Range("D8").Value = 0.02100840336134454
$ws.Range("F8").Value = 0.07352941176470588
$ws.Range("J8").Value = 0.1197478991596639
$ws.Range("O8").Value = 0.01470588235294118
$ws.Range("Q8").Value = 0.2647058823529412
$ws.Range("R8").Value = 0.07142857142857142
$ws.Range("S8").Value = 0.3445378151260504
$ws.Range("B9").Value = 0.07906976744186046
$ws.Range("D9").Value = 0.009302325581395349
$ws.Range("F9").Value = 0.02790697674418605
$ws.Range("J9").Value = 0.1116279069767442
$ws.Range("O9").Value = 0.0186046511627907
$ws.Range("Q9").Value = 0.2930232558139535
$ws.Range("R9").Value = 0.08372093023255814
$ws.Range("S9").Value = 0.3767441860465116
$ws.Range("B10").Value = 0.1180555555555556
$ws.Range("D10").Value = 0.02430555555555556
$ws.Range("E10").Value = 0.0008680555555555555
$ws.Range("F10").Value = 0.06510416666666667
$ws.Range("J10").Value = 0.09635416666666667
$ws.Range("O10").Value = 0.006076388888888889
$ws.Range("Q10").Value = 0.2951388888888889
$ws.Range("R10").Value = 0.07118055555555555
$ws.Range("S10").Value = 0.3229166666666667
$ws.Range("G11").Value = 0.1529850746268657
$ws.Range("J11").Value = 0.09701492537313433
$ws.Range("K11").Value = 0.2388059701492537
$ws.Range("L11").Value = 0.4925373134328358
$ws.Range("S11").Value = 0.01865671641791045
$ws.Range("G12").Value = 0.7714285714285715
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("K12").Value = 0.007142857142857143
$ws.Range("L12").Value = 0.05
$ws.Range("S12").Value = 0.02857142857142857
$ws.Range("G13").Value = 0.6896551724137931
$ws.Range("J13").Value = 0.3103448275862069
$ws.Range("F15").Value = 0.02380952380952381
$ws.Range("H15").Value = 0.2083333333333333
$ws.Range("I15").Value = 0.07738095238095238
$ws.Range("J15").Value = 0.2916666666666667
$ws.Range("K15").Value = 0.1071428571428571
$ws.Range("M15").Value = 0.01785714285714286
$ws.Range("O15").Value = 0.04166666666666666
$ws.Range("S15").Value = 0.2321428571428572
$ws.Range("F16").Value = 0.01522842639593909
$ws.Range("H16").Value = 0.2030456852791878
$ws.Range("I16").Value = 0.09137055837563451
$ws.Range("J16").Value = 0.3756345177664975
$ws.Range("K16").Value = 0.08121827411167512
$ws.Range("M16").Value = 0.04568527918781726
$ws.Range("O16").Value = 0.05583756345177665
$ws.Range("S16").Value = 0.1319796954314721
$ws.Range("F17").Value = 0.02791461412151067
$ws.Range("H17").Value = 0.187192118226601
$ws.Range("I17").Value = 0.1165845648604269
$ws.Range("J17").Value = 0.4187192118226601
$ws.Range("K17").Value = 0.07717569786535304
$ws.Range("M17").Value = 0.02298850574712644
$ws.Range("N17").Value = 0.003284072249589491
$ws.Range("O17").Value = 0.0361247947454844
$ws.Range("S17").Value = 0.110016420361248
$ws.Range("F18").Value = 0.01851851851851852
$ws.Range("H18").Value = 0.2098765432098765
$ws.Range("I18").Value = 0.07407407407407407
$ws.Range("J18").Value = 0.4135802469135803
$ws.Range("K18").Value = 0.1234567901234568
$ws.Range("M18").Value = 0.006172839506172839
$ws.Range("O18").Value = 0.04938271604938271
$ws.Range("S18").Value = 0.1049382716049383
$ws.Range("F19").Value = 0.01658986175115208
$ws.Range("H19").Value = 0.2350230414746544
$ws.Range("I19").Value = 0.09400921658986175
$ws.Range("J19").Value = 0.3594470046082949
$ws.Range("K19").Value = 0.08940092165898618
$ws.Range("M19").Value = 0.02949308755760369
$ws.Range("N19").Value = 0.0009216589861751152
$ws.Range("O19").Value = 0.06728110599078341
$ws.Range("S19").Value = 0.1078341013824885
